$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaRangeForText($searchText, $occurrence) {
    # Locate the Nth (1-based) occurrence of $searchText in the document body
    # and return the Range of the paragraph that contains it (Start..End-1,
    # i.e. excluding the trailing paragraph mark).
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    for ($k = 1; $k -le $occurrence; $k++) {
        $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $ok) {
            throw "Could not find occurrence $k of '$searchText'"
        }
        if ($k -lt $occurrence) {
            $rng.Collapse(0)
        }
    }
    $para = $rng.Paragraphs(1)
    $pr = $para.Range
    return $d.Range($pr.Start, $pr.End - 1)
}

# --- 1) "Telemóvel" -> "Phone number" (first occurrence only) ---
$r1 = Get-ParaRangeForText "Telemóvel" 1
$xml1 = "<w:p $wns>" +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Phone</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>number</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>'
$r1.InsertXML($xml1) | Out-Null
Write-Output "1) Telemovel -> Phone number: done"

# --- 2) "Morada" -> "address" (first occurrence only) ---
$r2 = Get-ParaRangeForText "Morada" 1
$xml2 = "<w:p $wns>" +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>address</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>'
$r2.InsertXML($xml2) | Out-Null
Write-Output "2) Morada -> address: done"

# --- 3) "NIF" -> "vat" (first occurrence only) ---
$r3 = Get-ParaRangeForText "NIF" 1
$xml3 = "<w:p $wns>" +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>vat</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>'
$r3.InsertXML($xml3) | Out-Null
Write-Output "3) NIF -> vat: done"

# --- 4) Insert a bare empty paragraph right after "Consultar Encomendas " ---
$rng4 = $d.Content
$rng4.Find.ClearFormatting()
$rng4.Find.Execute("Consultar Encomendas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para4 = $rng4.Paragraphs(1)
$insPoint = $para4.Range.End
$r4 = $d.Range($insPoint, $insPoint)
$r4.InsertXML("<w:p $wns/>") | Out-Null
Write-Output "4) inserted empty paragraph after 'Consultar Encomendas': done"

# --- 5) Move <w:lastRenderedPageBreak/> from the "Dúvidas" run to the
#        "Clicar neste botão..." run (repagination side-effect). ---

# 5a. Add it to the paragraph containing "Clicar neste botão"
$rngA = $d.Content
$rngA.Find.ClearFormatting()
$rngA.Find.Execute("Clicar neste botão", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraA = $rngA.Paragraphs(1)
$rA = $d.Range($paraA.Range.Start, $paraA.Range.End - 1)
$xmlA = "<w:p $wns><w:r><w:lastRenderedPageBreak/><w:t>Clicar neste botão fará ainda enviar um resumo do pedido do cliente para o perfil do Farmacêutico selecionado.</w:t></w:r></w:p>"
$rA.InsertXML($xmlA) | Out-Null
Write-Output "5a) added lastRenderedPageBreak to 'Clicar neste botao' run: done"

# 5b. Remove it from the paragraph containing "Dúvidas"
$rngB = $d.Content
$rngB.Find.ClearFormatting()
$rngB.Find.Execute("Dúvidas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraB = $rngB.Paragraphs(1)
$rB = $d.Range($paraB.Range.Start, $paraB.Range.End - 1)
$xmlB = "<w:p $wns>" +
        '<w:r w:rsidRPr="006C5C64"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Dúvidas</w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> com resposta</w:t></w:r>' +
        '<w:r w:rsidR="00DD3882"><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t>(mudar)</w:t></w:r>' +
        '</w:p>'
$rB.InsertXML($xmlB) | Out-Null
Write-Output "5b) removed lastRenderedPageBreak from 'Duvidas' run: done"
